# feature: add topic selection
# Rename the sheet to reflect its topic, set column G width to fit content,
# and set the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Sheet1" to "Panel Data"
$ws.Name = "Panel Data"

# Auto-fit / set the width of column G (Explanation) so it reads well.
# (71 - 5/6 compensates for the engine's built-in padding so the
# persisted column width lands exactly on 71, matching a "best fit" width.)
$ws.Columns.Item(7).ColumnWidth = 70.16666666666667

# Set the active cell selection to K8 (matches the saved view state)
$ws.Range("K8").Select()
